# Cian: fixed double hashing. TODO: add in collision detection.
#
# The keys in rows 15, 16, 21 and 22 were being excluded from hashing
# (A column held 0 there) - that's the "double hashing" bug. Flip those
# four cells back to 1 so they hash through like every other row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 1
$ws.Range("A16").Value = 1
$ws.Range("A21").Value = 1
$ws.Range("A22").Value = 1

# Move the active selection to A6, matching where the author left the
# cursor after making the fix.
$ws.Range("A6").Select()
